$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "name"
$ws.Range("A7").Value = "name"

$ws.Range("A7").Select()
